# Update (Removed Auto Arima)
# Fill in the previously-blank "Amazon Mean Forecast" column (D) and adjust
# the corresponding P70/P80/P90 forecast columns (E/F/G) on the
# "Forecast Comparison" sheet for rows 2-17 (weeks W01-W16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @(
    @{ Row = 2;  D = 35; E = 43; F = 52; G = 66 }
    @{ Row = 3;  D = 33; E = 40; F = 52; G = 71 }
    @{ Row = 4;  D = 26; E = 32; F = 39; G = 51 }
    @{ Row = 5;  D = 24; E = 30; F = 37; G = 48 }
    @{ Row = 6;  D = 22; E = 27; F = 36; G = 50 }
    @{ Row = 7;  D = 25; E = 30; F = 40; G = 57 }
    @{ Row = 8;  D = 24; E = 29; F = 40; G = 58 }
    @{ Row = 9;  D = 24; E = 29; F = 40; G = 60 }
    @{ Row = 10; D = 23; E = 27; F = 37; G = 54 }
    @{ Row = 11; D = 23; E = 28; F = 39; G = 58 }
    @{ Row = 12; D = 24; E = 28; F = 39; G = 59 }
    @{ Row = 13; D = 27; E = 32; F = 44; G = 65 }
    @{ Row = 14; D = 26; E = 31; F = 42; G = 63 }
    @{ Row = 15; D = 23; E = 27; F = 39; G = 60 }
    @{ Row = 16; D = 25; E = 30; F = 42; G = 63 }
    @{ Row = 17; D = 23; E = 27; F = 39; G = 60 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}
